$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each crypto row.
# Price cells are forced to remain text (NumberFormat "@") so values
# like "1.000" / "304.89" are not auto-converted to numbers by Excel,
# then the cell style is reset back to "Normal" so no stray number
# formatting is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.904.98"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.862.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5050"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.36%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3645"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07168"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8942"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07479"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.851.56"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.69"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.229"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008493"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.941.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.033"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.076.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.46%  "

$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.393"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.23%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("E26").Value = "  -2.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.078"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.702"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.674"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09238"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05097"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.72%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7516"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.994"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.95%  "

$ws.Range("E36").Value = "  -0.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.268"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.534"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5580"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.23%  "

$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.515"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.517"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1470"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4694"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("E48").Value = "  +0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.562"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.29%  "
